$wb = $excel.ActiveWorkbook

# The three sheets (illnessratio, illnessday, chronicratio) each have an
# identical 4-row descriptive-statistics table (Variable / Observations /
# Mean / Stdev / Min / Q1 / Median / Q3 / Max) for the rows:
#   row2: <sheet-specific dependent variable>
#   row3: income
#   row4: edu
#
# This edit re-estimates the model with province fixed effects entered as
# individual dummies (manually constructed), which changes the "income"
# row's coefficients (they are now on a different scale) and introduces a
# new "income2" variable row between "income" and "edu".

for ($i = 1; $i -le 3; $i++) {
    $ws = $wb.Worksheets.Item($i)

    # Updated coefficients/stats for the existing "income" row (row 3).
    $ws.Range("C3").Value = -0.0000000000000000129996609418723
    $ws.Range("D3").Value = 2.56148082142689
    $ws.Range("E3").Value = -4.98028433450118
    $ws.Range("F3").Value = -1.57457427074231
    $ws.Range("G3").Value = -0.636455752008566
    $ws.Range("H3").Value = 1.13768013087621
    $ws.Range("I3").Value = 11.0699229131785

    # Insert a brand-new row for "income2" right after "income" (row 4),
    # pushing the existing "edu" row down to row 5 untouched.
    $ws.Rows.Item(4).Insert()

    $ws.Range("A4").Value = "income2"
    $ws.Range("B4").Value = 282
    $ws.Range("C4").Value = 6.53791738861388
    $ws.Range("D4").Value = 15.1350572466699
    $ws.Range("E4").Value = 0.000406504120631675
    $ws.Range("F4").Value = 0.605673881433441
    $ws.Range("G4").Value = 1.97130056034504
    $ws.Range("H4").Value = 5.8019088939939
    $ws.Range("I4").Value = 122.543193303713
}
